$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update measurement values in rows 11-13 (columns C, D, E)
$ws.Range("C11").Value = 233
$ws.Range("D11").Value = 207
$ws.Range("E11").Value = 60

$ws.Range("C12").Value = 26
$ws.Range("D12").Value = 89
$ws.Range("E12").Value = 385

$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 206
$ws.Range("E13").Value = 294

# Update the selected / active cell to C14
$ws.Range("C14").Select()
